$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IJN_Destroyers")
$ws.Columns("C").Insert()
$ws.Columns("C").ColumnWidth = 42
$ws.Range("D22").Value = "Total"
$ws.Range("E22").Formula = "=SUM(E5:E21)"
